$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell ref, new value, whether it needs forced text format
# (D-column price values look numeric and Excel would silently coerce them
#  to Double and reformat/round the display text, so those are forced to
#  Text ("@") before the assignment, then the style is reset to Normal so
#  no stray number-format style is left attached to the cell.)
$updates = @(
    @{ Cell = 'D2'; Value = '64.529.80'; ForceText = $true }
    @{ Cell = 'E2'; Value = '  +1.05%  '; ForceText = $true }
    @{ Cell = 'D3'; Value = '3.170.94'; ForceText = $true }
    @{ Cell = 'E3'; Value = '  +0.14%  '; ForceText = $true }
    @{ Cell = 'E4'; Value = '  -0.04%  '; ForceText = $true }
    @{ Cell = 'D5'; Value = '572.13'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  +0.71%  '; ForceText = $true }
    @{ Cell = 'D6'; Value = '164.30'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  -2.01%  '; ForceText = $true }
    @{ Cell = 'E7'; Value = '  +0.00%  '; ForceText = $true }
    @{ Cell = 'D8'; Value = '0.582'; ForceText = $true }
    @{ Cell = 'E8'; Value = '  -4.73%  '; ForceText = $true }
    @{ Cell = 'E9'; Value = '  -2.30%  '; ForceText = $true }
    @{ Cell = 'E10'; Value = '  -0.88%  '; ForceText = $true }
    @{ Cell = 'E11'; Value = '  -0.07%  '; ForceText = $true }
    @{ Cell = 'D12'; Value = '3.727.41'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  +0.04%  '; ForceText = $true }
    @{ Cell = 'E13'; Value = '  -1.03%  '; ForceText = $true }
    @{ Cell = 'D14'; Value = '64.556.32'; ForceText = $true }
    @{ Cell = 'E14'; Value = '  +0.97%  '; ForceText = $true }
    @{ Cell = 'D15'; Value = '25.34'; ForceText = $true }
    @{ Cell = 'E15'; Value = '  +0.30%  '; ForceText = $true }
    @{ Cell = 'D16'; Value = '3.161.19'; ForceText = $true }
    @{ Cell = 'E16'; Value = '  -0.19%  '; ForceText = $true }
    @{ Cell = 'E17'; Value = '  -1.87%  '; ForceText = $true }
    @{ Cell = 'D18'; Value = '408.03'; ForceText = $true }
    @{ Cell = 'E18'; Value = '  -1.43%  '; ForceText = $true }
    @{ Cell = 'E19'; Value = '  +0.08%  '; ForceText = $true }
    @{ Cell = 'E20'; Value = '  -1.14%  '; ForceText = $true }
    @{ Cell = 'D21'; Value = '7.14'; ForceText = $true }
    @{ Cell = 'E21'; Value = '  +0.59%  '; ForceText = $true }
    @{ Cell = 'E22'; Value = '  +0.26%  '; ForceText = $true }
    @{ Cell = 'D23'; Value = '68.81'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  -2.72%  '; ForceText = $true }
    @{ Cell = 'D24'; Value = '0.488'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  -0.40%  '; ForceText = $true }
    @{ Cell = 'D25'; Value = '0.195'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  -4.02%  '; ForceText = $true }
    @{ Cell = 'D26'; Value = '0.0000103'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  -3.55%  '; ForceText = $true }
    @{ Cell = 'D27'; Value = '8.86'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  +1.89%  '; ForceText = $true }
    @{ Cell = 'D28'; Value = '0.992'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  -0.70%  '; ForceText = $true }
    @{ Cell = 'E29'; Value = '  -0.90%  '; ForceText = $true }
    @{ Cell = 'D30'; Value = '21.26'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  -2.25%  '; ForceText = $true }
    @{ Cell = 'E31'; Value = '  +0.17%  '; ForceText = $true }
    @{ Cell = 'D32'; Value = '4.89'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  -1.78%  '; ForceText = $true }
    @{ Cell = 'E33'; Value = '  -0.07%  '; ForceText = $true }
    @{ Cell = 'D34'; Value = '156.60'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  +0.48%  '; ForceText = $true }
    @{ Cell = 'E35'; Value = '  -1.39%  '; ForceText = $true }
    @{ Cell = 'E36'; Value = '  +0.09%  '; ForceText = $true }
    @{ Cell = 'D37'; Value = '2.693.67'; ForceText = $true }
    @{ Cell = 'E37'; Value = '  -1.52%  '; ForceText = $true }
    @{ Cell = 'D38'; Value = '24.10'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  -3.18%  '; ForceText = $true }
    @{ Cell = 'D39'; Value = '4.11'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  -1.19%  '; ForceText = $true }
    @{ Cell = 'D40'; Value = '0.696'; ForceText = $true }
    @{ Cell = 'E40'; Value = '  -2.25%  '; ForceText = $true }
    @{ Cell = 'D41'; Value = '0.0621'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  -0.41%  '; ForceText = $true }
    @{ Cell = 'D42'; Value = '5.49'; ForceText = $true }
    @{ Cell = 'E42'; Value = '  -2.60%  '; ForceText = $true }
    @{ Cell = 'B43'; Value = 'InjectiveProtocol'; ForceText = $false }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; ForceText = $false }
    @{ Cell = 'D43'; Value = '21.52'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  -1.15%  '; ForceText = $true }
    @{ Cell = 'B44'; Value = 'VeChain'; ForceText = $false }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; ForceText = $false }
    @{ Cell = 'D44'; Value = '0.0258'; ForceText = $true }
    @{ Cell = 'E44'; Value = '  -1.45%  '; ForceText = $true }
    @{ Cell = 'B45'; Value = 'Bittensor'; ForceText = $false }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; ForceText = $false }
    @{ Cell = 'D45'; Value = '291.30'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  -1.65%  '; ForceText = $true }
    @{ Cell = 'E46'; Value = '  +0.04%  '; ForceText = $true }
    @{ Cell = 'D47'; Value = '0.0985'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  -0.59%  '; ForceText = $true }
    @{ Cell = 'E48'; Value = '  -6.23%  '; ForceText = $true }
    @{ Cell = 'D49'; Value = '10.47'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  +0.26%  '; ForceText = $true }
    @{ Cell = 'E50'; Value = '  -0.74%  '; ForceText = $true }
    @{ Cell = 'D51'; Value = '0.879'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  -6.40%  '; ForceText = $true }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = '@'
        $rng.Value = $u.Value
        $rng.Style = 'Normal'
    } else {
        $rng.Value = $u.Value
    }
}
